$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all existing data rows
# (rows 2-369) from 45189 (2023-09-20) to 45190 (2023-09-21).
$ws.Range("C2:C369").Value2 = 45190

# Row 369 picks up an explicit row height (matches the default height of 15
# that every other data row already carries explicitly).
$ws.Rows.Item(369).RowHeight = 15

# Append the new record as row 370.
$ws.Cells.Item(370, 1).Value2 = "A 44441-2023"

$ws.Cells.Item(370, 2).Value2 = 45189
$ws.Cells.Item(370, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(370, 3).Value2 = 45190
$ws.Cells.Item(370, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(370, 4).Value2 = "VÄSTRA GÖTALANDS LÄN"
$ws.Cells.Item(370, 5).Value2 = "TRANEMO"

$ws.Cells.Item(370, 7).Value2 = 5.7
$ws.Cells.Item(370, 8).Value2 = 0
$ws.Cells.Item(370, 9).Value2 = 0
$ws.Cells.Item(370, 10).Value2 = 0
$ws.Cells.Item(370, 11).Value2 = 0
$ws.Cells.Item(370, 12).Value2 = 0
$ws.Cells.Item(370, 13).Value2 = 0
$ws.Cells.Item(370, 14).Value2 = 0
$ws.Cells.Item(370, 15).Value2 = 0
$ws.Cells.Item(370, 16).Value2 = 0
$ws.Cells.Item(370, 17).Value2 = 0

$ws.Cells.Item(370, 18).WrapText = $true
